# Insert a new (empty) "next page" section-break paragraph immediately
# before the paragraph that holds the "Sales Distribution" picture, and
# mark the document's final section as continuous.

$d = $word.ActiveDocument

# Locate the paragraph that contains the inline picture (it directly
# follows the "Visual comparison of unit sales across product lines"
# paragraph).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $target = $p
        break
    }
}

# Insert a section break (next page -> wdSectionBreakNextPage = 2) right
# at the start of that paragraph; this mints a new, empty paragraph just
# before it whose pPr carries the split-off sectPr.
$insertionPoint = $d.Range($target.Range.Start, $target.Range.Start)
$insertionPoint.InsertBreak(2)

# The document's last (trailing) section now needs to be continuous.
$lastSection = $d.Sections.Last
$lastSection.PageSetup.SectionStart = 0

Write-Output "done"
